$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.754.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.749.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.46%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.97%  '

$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5051'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.62%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.63'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2642'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +10.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06174'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.759.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06917'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6071'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.469'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.785.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006666'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.975.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.040'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.96%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.200'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.157'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.453'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.778'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08262'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.703'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.389'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04360'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9995'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.650'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9978'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5991'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.694'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.941'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01548'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.37%  '

$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.09'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7515'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3796'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.856'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05496'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1082'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.39%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.895'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.49%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.65%  '
